$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 50")

$ws.Range("C12").Value = 0.50347222222222221
$ws.Range("D12").Value = 0.51388888888888895
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = "De buttons in een List<Image> toegevoegd"
